$wb = $excel.ActiveWorkbook

$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("D2").Value = "2016-01-08 10:19:59"
$wsZh.Range("G2").Value = "2016-01-08 10:20:42"

$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("D2").Value = "2016-01-08 10:20:10"
$wsDe.Range("G2").Value = "2016-01-08 10:21:00"
